$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grade the "Generic" section (rows 3-6): fill in the "Points for grading"
# column (E) to match the max points already recorded in column D - the
# grader finished evaluating these rubric items with full marks.
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Grade the "Customer Class" section (rows 10-14) the same way.
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Leave the selection where the grader ended up after entering the last
# grade (matches the saved cursor position / scroll state).
$ws.Range("E15").Select()
